$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.912.76'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '3.520.72'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.52%  '
$ws.Range("D7").Value = '3.519.14'
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.491'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("E10").Value = '  +0.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.67%  '
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("D13").Value = '4.116.95'
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000181'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D17").Value = '3.516.67'
$ws.Range("E17").Value = '  -0.90%  '
$ws.Range("D18").Value = '64.938.09'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '390.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.575'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '3.660.07'
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  -2.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +15.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("E32").Value = '  +1.42%  '
$ws.Range("D33").Value = '3.522.88'
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.20'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.57%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("E37").Value = '  +1.26%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.12'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.12%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '168.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0817'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.823'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.83%  '
$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.93%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("E48").Value = '  -1.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.89'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").Value = '2.382.76'
$ws.Range("E50").Value = '  -2.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0268'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.30%  '
